$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '58.581.30'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  +1.38%  '

$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '2.517.60'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  +2.72%  '

$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E4').Value = '  -0.15%  '

$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '521.30'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +0.80%  '

$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '132.50'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +0.27%  '

$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.998'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  -0.15%  '

$ws.Range('E8').Value = '  +0.46%  '

$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '2.517.10'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +2.37%  '

$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.0974'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -0.41%  '

$ws.Range('E11').Value = '  -0.19%  '

$ws.Range('E12').Value = '  -2.16%  '

$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.332'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -1.83%  '

$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '2.952.81'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +2.25%  '

$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '58.525.69'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +1.32%  '

$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '22.09'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -0.40%  '

$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '2.507.76'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  +2.27%  '

$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '10.65'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +0.41%  '

$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '322.07'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +0.85%  '

$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '4.16'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +0.13%  '

$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '6.07'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +6.58%  '

$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -0.01%  '

$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '64.48'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +0.22%  '

$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '0.400'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -1.77%  '

$ws.Range('E26').Value = '  -0.03%  '

$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '0.160'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +0.64%  '

$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '7.37'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +0.55%  '

$ws.Range('D29').Value = '0.0₃0754'
$ws.Range('E29').Value = '  +2.16%  '

$ws.Range('E30').Value = '  +2.31%  '

$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '167.90'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +0.09%  '

$ws.Range('E32').Value = '  +3.52%  '

$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '6.25'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +0.77%  '

$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -0.10%  '

$ws.Range('E36').Value = '  +0.42%  '

$ws.Range('E37').Value = '  -2.93%  '

$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '3.96'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +0.30%  '

$ws.Range('E39').Value = '  +0.34%  '

$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '35.94'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -0.67%  '

$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.778'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -0.83%  '

$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '3.49'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +1.94%  '

$ws.Range('B43').Value = 'Bittensor'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '277.80'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +2.96%  '

$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '5.02'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +1.47%  '

$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.598'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +2.20%  '

$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '123.53'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -0.67%  '

$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.0918'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +1.09%  '

$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '0.0501'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +3.27%  '

$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '17.83'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +1.79%  '

$ws.Range('E50').Value = '  +1.84%  '

$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '16.91'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +1.17%  '
